# Generate Report for Handoff
# Updates the localization-status report: the "In Translation" rows have
# moved to "Ready for handoff", with refreshed handoff timestamps on the
# Overview sheet and each language sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-17 10:36:54"
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 16.3333333333333

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-17 10:36:49"
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 16.3333333333333

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-17 10:36:54"
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 16.3333333333333
